# Auto-generated: apply weekly crime-data update (2022-08-25) to violent-crime workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("I2").Value = 4666
$ws.Range("I3").Value = 4853
$ws.Range("I4").Value = 1106
$ws.Range("I6").Value = 5293
$ws.Range("I7").Value = 16362

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("I3").Value = 48
$ws.Range("I4").Value = 13

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("I2").Value = 168
$ws.Range("I3").Value = 166
$ws.Range("I6").Value = 139
$ws.Range("I7").Value = 520

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("I3").Value = 114
$ws.Range("I7").Value = 309

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("I2").Value = 153
$ws.Range("I6").Value = 210
$ws.Range("I7").Value = 646

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("I6").Value = 57
$ws.Range("I7").Value = 156

$ws = $wb.Worksheets.Item("New City")
$ws.Range("I3").Value = 109
$ws.Range("I6").Value = 107
$ws.Range("I7").Value = 367

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("I2").Value = 129
$ws.Range("I4").Value = 61
$ws.Range("I6").Value = 112
$ws.Range("I7").Value = 523
$ws.Range("I8").Value = 998
$ws.Range("I11").Value = 247
$ws.Range("I16").Value = 43
$ws.Range("I18").Value = 113
$ws.Range("I19").Value = 461
$ws.Range("I20").Value = 398
$ws.Range("I22").Value = 44
$ws.Range("I29").Value = 1040
$ws.Range("I31").Value = 156
$ws.Range("I33").Value = 757
$ws.Range("I36").Value = 222
$ws.Range("I37").Value = 520
$ws.Range("I41").Value = 72
$ws.Range("I42").Value = 556
$ws.Range("I43").Value = 131
$ws.Range("I46").Value = 32
$ws.Range("I47").Value = 112
$ws.Range("I50").Value = 74
$ws.Range("I54").Value = 363
$ws.Range("I55").Value = 179
$ws.Range("I57").Value = 61
$ws.Range("I60").Value = 85
$ws.Range("I63").Value = 60
$ws.Range("I65").Value = 367
$ws.Range("I66").Value = 46
$ws.Range("I67").Value = 646
$ws.Range("I73").Value = 137
$ws.Range("I77").Value = 100
$ws.Range("I79").Value = 457
$ws.Range("I83").Value = 340
$ws.Range("I85").Value = 736
$ws.Range("I86").Value = 96
$ws.Range("I90").Value = 204
$ws.Range("I92").Value = 49
$ws.Range("I93").Value = 97
$ws.Range("I95").Value = 269
$ws.Range("I99").Value = 309
$ws.Range("I101").Value = 16362

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("I2").Value = 125
$ws.Range("I4").Value = 12
$ws.Range("I7").Value = 340

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("I3").Value = 103
$ws.Range("I7").Value = 269

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("I2").Value = 174
$ws.Range("I4").Value = 33
$ws.Range("I6").Value = 239
$ws.Range("I7").Value = 757

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("I2").Value = 82
$ws.Range("I3").Value = 77
$ws.Range("I6").Value = 176
$ws.Range("I7").Value = 363

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("I2").Value = 303
$ws.Range("I3").Value = 361
$ws.Range("I6").Value = 284
$ws.Range("I7").Value = 1040

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("I2").Value = 170
$ws.Range("I6").Value = 127
$ws.Range("I7").Value = 461

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("I2").Value = 193
$ws.Range("I3").Value = 295
$ws.Range("I7").Value = 736

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("I3").Value = 31
$ws.Range("I6").Value = 26
$ws.Range("I7").Value = 112

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("I2").Value = 23
$ws.Range("I7").Value = 72

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("I3").Value = 188
$ws.Range("I4").Value = 44
$ws.Range("I7").Value = 556

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("I6").Value = 56
$ws.Range("I7").Value = 179

$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Range("I2").Value = 7
$ws.Range("I7").Value = 32

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("I3").Value = 148
$ws.Range("I7").Value = 457

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("I2").Value = 110
$ws.Range("I6").Value = 127
$ws.Range("I7").Value = 398

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("I2").Value = 33
$ws.Range("I7").Value = 113

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("I3").Value = 71
$ws.Range("I7").Value = 222

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("I6").Value = 40
$ws.Range("I7").Value = 97

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("I2").Value = 23
$ws.Range("I7").Value = 112

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("I2").Value = 19
$ws.Range("I7").Value = 74

$ws = $wb.Worksheets.Item("North Center")
$ws.Range("I6").Value = 20
$ws.Range("I7").Value = 46

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("I2").Value = 110
$ws.Range("I7").Value = 247

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("I3").Value = 45
$ws.Range("I7").Value = 137

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("I2").Value = 43
$ws.Range("I7").Value = 129

$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("I6").Value = 21
$ws.Range("I7").Value = 49

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("I2").Value = 318
$ws.Range("I6").Value = 317
$ws.Range("I7").Value = 998

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("I6").Value = 23
$ws.Range("I7").Value = 96

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("I3").Value = 45
$ws.Range("I7").Value = 204

$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("I3").Value = 12
$ws.Range("I7").Value = 61

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("I6").Value = 25
$ws.Range("I7").Value = 85

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("I6").Value = 74
$ws.Range("I7").Value = 131

$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("I6").Value = 13
$ws.Range("I7").Value = 44

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("I3").Value = 35
$ws.Range("I7").Value = 100

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("I3").Value = 163
$ws.Range("I7").Value = 523

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("I6").Value = 17
$ws.Range("I7").Value = 61

$ws = $wb.Worksheets.Item("Bucktown")
$ws.Range("I2").Value = 9
$ws.Range("I7").Value = 43
